$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "282.37"
Set-TextValue $ws.Range("E2") "1.93%"
Set-TextValue $ws.Range("E3") "1.78%"
Set-TextValue $ws.Range("D4") "4.870"
Set-TextValue $ws.Range("E4") "-0.37%"
Set-TextValue $ws.Range("D5") "0.06507"
Set-TextValue $ws.Range("E5") "1.32%"
Set-TextValue $ws.Range("D6") "7.119"
Set-TextValue $ws.Range("E6") "2.49%"
Set-TextValue $ws.Range("D7") "1.279"
Set-TextValue $ws.Range("E7") "3.03%"
Set-TextValue $ws.Range("D8") "0.9170"
Set-TextValue $ws.Range("E8") "3.97%"
Set-TextValue $ws.Range("D9") "0.1552"
Set-TextValue $ws.Range("E9") "2.21%"
Set-TextValue $ws.Range("D10") "0.06447"
Set-TextValue $ws.Range("E10") "28.02%"
Set-TextValue $ws.Range("D11") "0.07543"
Set-TextValue $ws.Range("E11") "0.62%"
Set-TextValue $ws.Range("D12") "0.02923"
Set-TextValue $ws.Range("E12") "0.37%"
Set-TextValue $ws.Range("D13") "0.08964"
Set-TextValue $ws.Range("E13") "-0.34%"
Set-TextValue $ws.Range("D14") "0.001587"
Set-TextValue $ws.Range("E14") "1.38%"
Set-TextValue $ws.Range("D15") "0.0006458"
Set-TextValue $ws.Range("E15") "0.83%"
Set-TextValue $ws.Range("D16") "0.006011"
Set-TextValue $ws.Range("E16") "5.53%"
Set-TextValue $ws.Range("D17") "3.503"
Set-TextValue $ws.Range("E17") "1.21%"
Set-TextValue $ws.Range("D18") "3.337"
Set-TextValue $ws.Range("E18") "0.62%"
Set-TextValue $ws.Range("D19") "2.227"
Set-TextValue $ws.Range("E19") "-1.96%"
Set-TextValue $ws.Range("D20") "0.3146"
Set-TextValue $ws.Range("E20") "0.31%"
Set-TextValue $ws.Range("D21") "0.1350"
Set-TextValue $ws.Range("E21") "0.94%"
Set-TextValue $ws.Range("D22") "3.874"
Set-TextValue $ws.Range("E22") "-1.01%"
Set-TextValue $ws.Range("D23") "0.04366"
Set-TextValue $ws.Range("E23") "-1.48%"
Set-TextValue $ws.Range("D25") "0.001164"
Set-TextValue $ws.Range("E25") "-0.84%"
Set-TextValue $ws.Range("D26") "0.004380"
Set-TextValue $ws.Range("E26") "13.45%"
Set-TextValue $ws.Range("D28") "0.0001171"
Set-TextValue $ws.Range("E28") "-2.42%"
Set-TextValue $ws.Range("D29") "0.0001651"
Set-TextValue $ws.Range("E29") "-14.72%"
Set-TextValue $ws.Range("D40") "0.04118"
Set-TextValue $ws.Range("E40") "-0.68%"
Set-TextValue $ws.Range("D41") "0.1395"
Set-TextValue $ws.Range("E41") "18.61%"
Set-TextValue $ws.Range("D42") "0.006278"
Set-TextValue $ws.Range("E42") "-7.61%"
Set-TextValue $ws.Range("D43") "0.002075"
Set-TextValue $ws.Range("E43") "-13.28%"
Set-TextValue $ws.Range("D44") "0.01166"
Set-TextValue $ws.Range("E44") "-0.47%"
Set-TextValue $ws.Range("D45") "0.00005528"
Set-TextValue $ws.Range("E45") "6.34%"
Set-TextValue $ws.Range("D46") "1.561"
Set-TextValue $ws.Range("E46") "5.01%"
Set-TextValue $ws.Range("D47") "0.01842"
Set-TextValue $ws.Range("E47") "-9.02%"
